$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.032.76"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.830.09"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.58%  "

$ws.Range("D5").Value = "'311.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").Value = "'0.4633"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.3709"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.09%  "

$ws.Range("D9").Value = "'0.07348"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").Value = "'0.8786"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.59%  "

$ws.Range("D11").Value = "'0.07888"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.22%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "1.780.20"
$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").Value = "'5.342"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "'6.540"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "'0.000008857"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.78%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").Value = "'14.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("D21").Value = "27.050.29"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").Value = "'5.109"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").Value = "2.000.50"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("D25").Value = "'152.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("D26").Value = "'1.846"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("E27").Value = "  +1.52%  "

$ws.Range("D28").Value = "'2.046"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.14%  "

$ws.Range("D29").Value = "'5.127"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").Value = "'115.84"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").Value = "'0.08892"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").Value = "'2.959"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").Value = "'0.7283"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "'4.440"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("D35").Value = "'1.134"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "'2.464"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").Value = "'1.077"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("E38").Value = "  +2.24%  "

$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").Value = "'2.954"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").Value = "'0.5173"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.43%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'8.182"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "'0.4842"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'1.005"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").Value = "'10.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").Value = "'102.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("D49").Value = "'1.631"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").Value = "  -0.51%  "

$ws.Range("E51").Value = "  +0.62%  "
